$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Volume/number header text: "43" -> "44"
$ws.Range("A8").Value = "Volume 30   Number  44"

# Reporting week date range update
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# Weekly crime-complaints table (rows 14-30): new collected figures
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -16.666666666666
$ws.Range("I14").Value = 110
$ws.Range("J14").Value = 111
$ws.Range("K14").Value = -0.9009009009
$ws.Range("L14").Value = -14.0625
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -74.943052391799
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -71.428571428571
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 316
$ws.Range("J15").Value = 336
$ws.Range("K15").Value = -5.95238095238
$ws.Range("L15").Value = 0.636942675159
$ws.Range("M15").Value = 24.409448818897
$ws.Range("N15").Value = -49.358974358974
$ws.Range("C16").Value = 107
$ws.Range("D16").Value = 107
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 420
$ws.Range("G16").Value = 419
$ws.Range("H16").Value = 0.238663484486
$ws.Range("I16").Value = 4178
$ws.Range("J16").Value = 4343
$ws.Range("K16").Value = -3.799217131015
$ws.Range("L16").Value = 28.672620880813
$ws.Range("M16").Value = 8.350622406639
$ws.Range("N16").Value = -69.881776239907
$ws.Range("C17").Value = 158
$ws.Range("D17").Value = 171
$ws.Range("E17").Value = -7.602339181286
$ws.Range("F17").Value = 584
$ws.Range("G17").Value = 563
$ws.Range("H17").Value = 3.730017761989
$ws.Range("I17").Value = 6932
$ws.Range("J17").Value = 6291
$ws.Range("K17").Value = 10.189159116197
$ws.Range("L17").Value = 29.983123945246
$ws.Range("M17").Value = 82.373059721126
$ws.Range("N17").Value = -10.762100926879
$ws.Range("C18").Value = 53
$ws.Range("D18").Value = 46
$ws.Range("E18").Value = 15.217391304347
$ws.Range("F18").Value = 200
$ws.Range("G18").Value = 223
$ws.Range("H18").Value = -10.313901345291
$ws.Range("I18").Value = 2478
$ws.Range("J18").Value = 2509
$ws.Range("K18").Value = -1.235552012754
$ws.Range("L18").Value = 34.018388318009
$ws.Range("M18").Value = -11.971580817051
$ws.Range("N18").Value = -84.462975735155
$ws.Range("C19").Value = 144
$ws.Range("D19").Value = 144
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 639
$ws.Range("G19").Value = 654
$ws.Range("H19").Value = -2.293577981651
$ws.Range("I19").Value = 6797
$ws.Range("J19").Value = 6865
$ws.Range("K19").Value = -0.990531682447
$ws.Range("L19").Value = 20.237042278436
$ws.Range("M19").Value = 70.564617314931
$ws.Range("N19").Value = 5.102829751043
$ws.Range("C20").Value = 91
$ws.Range("D20").Value = 68
$ws.Range("E20").Value = 33.823529411764
$ws.Range("F20").Value = 353
$ws.Range("G20").Value = 323
$ws.Range("H20").Value = 9.287925696594
$ws.Range("I20").Value = 4456
$ws.Range("J20").Value = 3294
$ws.Range("K20").Value = 35.276259866423
$ws.Range("L20").Value = 79.895034315704
$ws.Range("M20").Value = 149.217002237136
$ws.Range("N20").Value = -65.997710797405
$ws.Range("C21").Value = 557
$ws.Range("D21").Value = 544
$ws.Range("E21").Value = 2.389705882352
$ws.Range("F21").Value = 2221
$ws.Range("G21").Value = 2218
$ws.Range("H21").Value = 0.135256988277
$ws.Range("I21").Value = 25267
$ws.Range("J21").Value = 23749
$ws.Range("K21").Value = 6.391848077813
$ws.Range("L21").Value = 32.977211725698
$ws.Range("M21").Value = 52.128364139924
$ws.Range("N21").Value = -56.603805990656
$ws.Range("C22").Value = 12
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 43
$ws.Range("G22").Value = 28
$ws.Range("H22").Value = 53.571428571428
$ws.Range("I22").Value = 271
$ws.Range("J22").Value = 313
$ws.Range("K22").Value = -13.418530351437
$ws.Range("L22").Value = 23.181818181818
$ws.Range("M22").Value = 1.498127340823
$ws.Range("C23").Value = 35
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = 9.375
$ws.Range("F23").Value = 127
$ws.Range("G23").Value = 107
$ws.Range("H23").Value = 18.691588785046
$ws.Range("I23").Value = 1501
$ws.Range("J23").Value = 1351
$ws.Range("K23").Value = 11.102886750555
$ws.Range("L23").Value = 42.409867172675
$ws.Range("M23").Value = 65.85635359116
$ws.Range("C24").Value = 280
$ws.Range("D24").Value = 339
$ws.Range("E24").Value = -17.40412979351
$ws.Range("F24").Value = 1343
$ws.Range("G24").Value = 1389
$ws.Range("H24").Value = -3.311735061195
$ws.Range("I24").Value = 15232
$ws.Range("J24").Value = 15769
$ws.Range("K24").Value = -3.40541568901
$ws.Range("L24").Value = 36.658891082002
$ws.Range("M24").Value = 38.91472868217
$ws.Range("C25").Value = 202
$ws.Range("D25").Value = 189
$ws.Range("E25").Value = 6.878306878306
$ws.Range("F25").Value = 822
$ws.Range("G25").Value = 746
$ws.Range("H25").Value = 10.187667560321
$ws.Range("I25").Value = 8933
$ws.Range("J25").Value = 8486
$ws.Range("K25").Value = 5.267499410794
$ws.Range("L25").Value = 20.895926377047
$ws.Range("M25").Value = -6.254591247769
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("G26").Value = 47
$ws.Range("H26").Value = -10.63829787234
$ws.Range("I26").Value = 546
$ws.Range("J26").Value = 572
$ws.Range("K26").Value = -4.545454545454
$ws.Range("L26").Value = 6.432748538011
$ws.Range("C27").Value = 21
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 88
$ws.Range("G27").Value = 78
$ws.Range("H27").Value = 12.820512820512
$ws.Range("I27").Value = 907
$ws.Range("J27").Value = 790
$ws.Range("K27").Value = 14.810126582278
$ws.Range("L27").Value = 15.101522842639
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = 31.25
$ws.Range("I28").Value = 340
$ws.Range("J28").Value = 416
$ws.Range("K28").Value = -18.26923076923
$ws.Range("L28").Value = -34.235976789168
$ws.Range("M28").Value = -17.675544794188
$ws.Range("N28").Value = -72.536348949919
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 20
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = 25
$ws.Range("I29").Value = 283
$ws.Range("J29").Value = 353
$ws.Range("K29").Value = -19.830028328611
$ws.Range("L29").Value = -35.091743119266
$ws.Range("M29").Value = -18.208092485549
$ws.Range("N29").Value = -74.573225516621
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = -40
$ws.Range("I30").Value = 21
$ws.Range("J30").Value = 39
$ws.Range("K30").Value = -46.153846153846
$ws.Range("L30").Value = -47.5
